$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 155; this shifts the existing rows 155-170
# down to 156-171 (including their formatting), matching the diff.
$ws.Rows(155).Insert()

# Populate the newly inserted row 155 with the new record. The non-numeric
# "template" columns (Mercado ID, Mercado, Region, Codreg, Categoria ID,
# Categoria, Variedad, Calidad, Unidad de comercializacion, Origen, Kg o
# Unidades, Clasificacion) are identical to the former row 155 (now row
# 156); only Fecha, Volumen, Precio minimo/maximo/promedio and Precio $/Kg
# change.
$ws.Cells.Item(155, 1).Value = 10
$ws.Cells.Item(155, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(155, 3).Value = "La Araucanía"
$ws.Cells.Item(155, 4).Value = 44578
$ws.Cells.Item(155, 5).Value = 9
$ws.Cells.Item(155, 6).Value = 100112005
$ws.Cells.Item(155, 7).Value = "Puerro"
$ws.Cells.Item(155, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 40
$ws.Cells.Item(155, 11).Value = 14000
$ws.Cells.Item(155, 12).Value = 14000
$ws.Cells.Item(155, 13).Value = 14000
$ws.Cells.Item(155, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(155, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(155, 16).Value = 1167
$ws.Cells.Item(155, 17).Value = 12
$ws.Cells.Item(155, 18).Value = "Hortaliza"

# Keep the date cell formatted like the rest of column D.
$ws.Cells.Item(155, 4).NumberFormat = $ws.Cells.Item(156, 4).NumberFormat
